$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row imported (sales team / lost reason opportunity info)
$ws.Range("A3").Value = "My opportunity"
$ws.Range("B3").Value = "Twenty 2TB hard Disk"

# Widen column A slightly to fit the new content
$ws.Columns.Item(1).ColumnWidth = 13.5

# Move selection to E4, matching the post-edit cursor position
$ws.Range("E4").Select() | Out-Null
